# Automatische test-sync: 2025-06-29 14:46:50
# Append a new test-mail log row (row 18) to the "Logs" sheet, extend the
# conditional-formatting ranges that covered rows 2-17 to now cover rows
# 2-18, and bump the "Productinformatie" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: add row 18 -------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Wat zijn de verzendkosten?"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Testmail #3: Wat zijn de verzendkosten?"
$logs.Range("D18").Value = "Productinformatie"
$logs.Range("F18").Value = "2025-06-29 14:46:04"
$logs.Range("G18").Value = "Nee"
$logs.Range("H18").Value = "Ja"
$logs.Range("I18").Value = "Nee"

# --- Extend conditional formatting ranges from row 17 to row 18 ------------
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))

# --- Dashboard sheet: bump Productinformatie count from 4 to 5 -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 5
